# Added 'hole_id' index column to the 'train' sheet so cross validation can
# be performed: column A, which previously held a bare numeric row index
# (0..36), now holds the borehole "hole_id" string and a "hole_id" header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# hole_id values, in row order (row 2 .. row 38)
$holeIds = @(
    "LBU_05_15", "LBU_07_02", "MHZ_12_04", "LBU_05_29", "LBU_05_19",
    "LBU_05_30", "LBU_05_23", "LBU_05_17", "LBU_05_09", "MHZ_08_01",
    "LBU_05_13", "MHZ_08_05", "LBU_02_4",  "LBU_05_24", "LBU_05_10",
    "LBU_05_06", "MHZ_08_04", "LBU_05_26", "LBU_05_27", "MHZ_08_03",
    "LBU_05_16", "LBU_05_14", "LBU_05_21", "LBU_05_02", "LBU_05_03",
    "MHZ_12_02", "LBU_01_2",  "LBU_05_12", "LBU_01_1",  "LBU_05_11",
    "LBU_05_01", "LBU_05_07", "LBU_05_05", "LBU_05_04", "LBU_05_08",
    "MHZ_12_01", "LBU_07_03"
)

# A1 was previously empty (header row started at B1). Give it the same
# header formatting (bold, bordered, centered) as the other header cells
# before writing the "hole_id" label.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "hole_id"

# Replace the bare numeric index (0..36) in A2:A38 with the hole_id text.
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
